$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Shift the September transaction log (columns R = Details, S = Date)
# down by one row, from the bottom up, to make room for a new entry
# at the top (row 39). Rows 39-139 move to rows 40-140.
for ($r = 139; $r -ge 39; $r--) {
    $srcRVal = $ws.Cells.Item($r, 18).Value2      # column R
    $srcSVal = $ws.Cells.Item($r, 19).Value2      # column S
    $ws.Cells.Item($r + 1, 18).Value = $srcRVal
    $ws.Cells.Item($r + 1, 19).Value = $srcSVal
}

# New latest entry at the top of the shifted range.
$ws.Cells.Item(39, 18).Value = "balance your axis"
$ws.Cells.Item(39, 19).Value = "2024-09-15 07:56:24"

# The "Broadband" entry in the Group column (A) moves from row 148
# down to the newly added row 149.
$ws.Range("A148").Value = ""
$ws.Range("A149").Value = "Broadband"
